$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old location.road data column and everything from G onward
# (details.*, advisories.message, responders.arrivalTime/role, ...);
# only advisories.type, responders.agency and responders.personnel.name
# survive, moved into D:F.
$ws.Range("C2:C9").ClearContents()
$ws.Range("G1:S9").ClearContents()

# Update header row
$ws.Range("C1").Value = "location.type"
$ws.Range("D1").Value = "advisories.type"
$ws.Range("E1").Value = "responders.personnel.name"
$ws.Range("F1").Value = "responders.agency"

# Row 2 - move N2->D2, P2->F2, R2->E2
$ws.Range("D2").Value = "Diversion"
$ws.Range("E2").Value = "Sgt. Tan Wei 1"
$ws.Range("F2").Value = "Traffic Police 1"

# Row 3 - move N3->D3, R3->E3
$ws.Range("D3").Value = "Congestion Alert"
$ws.Range("E3").Value = "Cpl. Lim Hui 2"

# Row 4 - move P4->F4, R4->E4
$ws.Range("E4").Value = "Lt. Ravi Kumar 3"
$ws.Range("F4").Value = "SCDF 2"

# Row 5 - move R5->E5
$ws.Range("E5").Value = "Spec. Ong Jia 4"

# Row 6
$ws.Range("D6").Value = "Diversion 1"
$ws.Range("E6").Value = "Sgt. Tan Wei 1"
$ws.Range("F6").Value = "SCDF"

# Row 7
$ws.Range("D7").Value = "Congestion Alert 2"
$ws.Range("E7").Value = "Cpl. Lim Hui 2"

# Row 8
$ws.Range("E8").Value = "Lt. Ravi Kumar 3"
$ws.Range("F8").Value = "Traffic Police"

# Row 9
$ws.Range("E9").Value = "Spec. Ong Jia 4"
